# Update market-price / profit figures across the Leve-flipping sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to the latest scrape snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(26, 8).Value = 50000
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 14).Value = $null
$ws.Cells.Item(33, 8).Value = 96.565216
$ws.Cells.Item(33, 9).Value = 91.2381
$ws.Cells.Item(33, 10).Value = 152.5
$ws.Cells.Item(33, 11).Value = 91.2381
$ws.Cells.Item(33, 12).Value = 152.5
$ws.Cells.Item(33, 13).Value = 137.7619
$ws.Cells.Item(33, 14).Value = -610.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 15779.785
$ws.Cells.Item(32, 9).Value = 2810.817
$ws.Cells.Item(32, 11).Value = 2810.817
$ws.Cells.Item(32, 13).Value = -2523.817
$ws.Cells.Item(102, 8).Value = 3620.6667
$ws.Cells.Item(102, 9).Value = 3937.3635
$ws.Cells.Item(102, 10).Value = 2749.75
$ws.Cells.Item(102, 11).Value = 3937.3635
$ws.Cells.Item(102, 12).Value = 2749.75
$ws.Cells.Item(102, 13).Value = -2315.3635
$ws.Cells.Item(102, 14).Value = -5993.75
$ws.Cells.Item(122, 8).Value = 22100
$ws.Cells.Item(122, 9).Value = 27000
$ws.Cells.Item(122, 10).Value = 2500
$ws.Cells.Item(122, 11).Value = 81000
$ws.Cells.Item(122, 12).Value = 7500
$ws.Cells.Item(122, 13).Value = -78550
$ws.Cells.Item(122, 14).Value = -12400

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2518.0908
$ws.Cells.Item(20, 9).Value = 2298.3333
$ws.Cells.Item(20, 11).Value = 2298.3333
$ws.Cells.Item(20, 13).Value = -2051.3333
$ws.Cells.Item(105, 8).Value = 255132.05
$ws.Cells.Item(105, 9).Value = 6101.037
$ws.Cells.Item(105, 10).Value = 772350.3
$ws.Cells.Item(105, 11).Value = 6101.037
$ws.Cells.Item(105, 12).Value = 772350.3
$ws.Cells.Item(105, 13).Value = -4354.037
$ws.Cells.Item(105, 14).Value = -775844.3
$ws.Cells.Item(134, 8).Value = 23258542
$ws.Cells.Item(134, 9).Value = 29414024
$ws.Cells.Item(134, 10).Value = 4502.8887
$ws.Cells.Item(134, 11).Value = 88242072
$ws.Cells.Item(134, 12).Value = 13508.6661
$ws.Cells.Item(134, 13).Value = -88239537
$ws.Cells.Item(134, 14).Value = -18578.6661

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 102160
$ws.Cells.Item(16, 9).Value = 251000
$ws.Cells.Item(16, 10).Value = 2933.3333
$ws.Cells.Item(16, 11).Value = 251000
$ws.Cells.Item(16, 12).Value = 2933.3333
$ws.Cells.Item(16, 13).Value = -250713
$ws.Cells.Item(16, 14).Value = -3507.3333
$ws.Cells.Item(29, 8).Value = 1200
$ws.Cells.Item(29, 10).Value = 1200
$ws.Cells.Item(29, 12).Value = 1200
$ws.Cells.Item(29, 14).Value = -1786
$ws.Cells.Item(31, 8).Value = 2244.2942
$ws.Cells.Item(31, 9).Value = 1388.5264
$ws.Cells.Item(31, 10).Value = 3328.2666
$ws.Cells.Item(31, 11).Value = 1388.5264
$ws.Cells.Item(31, 12).Value = 3328.2666
$ws.Cells.Item(31, 13).Value = -1093.5264
$ws.Cells.Item(31, 14).Value = -3918.2666
$ws.Cells.Item(34, 8).Value = 2244.2942
$ws.Cells.Item(34, 9).Value = 1388.5264
$ws.Cells.Item(34, 10).Value = 3328.2666
$ws.Cells.Item(34, 11).Value = 1388.5264
$ws.Cells.Item(34, 12).Value = 3328.2666
$ws.Cells.Item(34, 13).Value = -1186.5264
$ws.Cells.Item(34, 14).Value = -3732.2666
$ws.Cells.Item(45, 8).Value = 5067
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 14).Value = $null
$ws.Cells.Item(86, 8).Value = 17861004
$ws.Cells.Item(86, 9).Value = 23813900
$ws.Cells.Item(86, 10).Value = 2314.2856
$ws.Cells.Item(86, 11).Value = 23813900
$ws.Cells.Item(86, 12).Value = 2314.2856
$ws.Cells.Item(86, 13).Value = -23812777
$ws.Cells.Item(86, 14).Value = -4560.2856
$ws.Cells.Item(89, 8).Value = 17861004
$ws.Cells.Item(89, 9).Value = 23813900
$ws.Cells.Item(89, 10).Value = 2314.2856
$ws.Cells.Item(89, 11).Value = 119069500
$ws.Cells.Item(89, 12).Value = 11571.428
$ws.Cells.Item(89, 13).Value = -119063884
$ws.Cells.Item(89, 14).Value = -22803.428
$ws.Cells.Item(107, 8).Value = 380.3889
$ws.Cells.Item(107, 9).Value = 222.66667
$ws.Cells.Item(107, 10).Value = 538.1111
$ws.Cells.Item(107, 11).Value = 222.66667
$ws.Cells.Item(107, 12).Value = 538.1111
$ws.Cells.Item(107, 13).Value = 1697.33333
$ws.Cells.Item(107, 14).Value = -4378.1111
$ws.Cells.Item(113, 8).Value = 102160
$ws.Cells.Item(113, 9).Value = 251000
$ws.Cells.Item(113, 10).Value = 2933.3333
$ws.Cells.Item(113, 11).Value = 251000
$ws.Cells.Item(113, 12).Value = 2933.3333
$ws.Cells.Item(113, 13).Value = -248830
$ws.Cells.Item(113, 14).Value = -7273.3333
$ws.Cells.Item(132, 8).Value = 2166.087
$ws.Cells.Item(132, 9).Value = 1820
$ws.Cells.Item(132, 10).Value = 4473.3335
$ws.Cells.Item(132, 11).Value = 5460
$ws.Cells.Item(132, 12).Value = 13420.0005
$ws.Cells.Item(132, 13).Value = -2930
$ws.Cells.Item(132, 14).Value = -18480.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 2641.6956
$ws.Cells.Item(139, 9).Value = 2232.8823
$ws.Cells.Item(139, 10).Value = 3800
$ws.Cells.Item(139, 11).Value = 6698.646900000001
$ws.Cells.Item(139, 12).Value = 11400
$ws.Cells.Item(139, 13).Value = -1558.646900000001
$ws.Cells.Item(139, 14).Value = -21680

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(13, 8).Value = 3006
$ws.Cells.Item(13, 10).Value = 3006
$ws.Cells.Item(13, 12).Value = 3006
$ws.Cells.Item(13, 14).Value = -3284
$ws.Cells.Item(19, 8).Value = 32401
$ws.Cells.Item(19, 10).Value = 50000
$ws.Cells.Item(19, 12).Value = 50000
$ws.Cells.Item(19, 14).Value = -50576
$ws.Cells.Item(29, 8).Value = 7469.3335
$ws.Cells.Item(29, 10).Value = 7469.3335
$ws.Cells.Item(29, 12).Value = 7469.3335
$ws.Cells.Item(29, 14).Value = -8049.3335
$ws.Cells.Item(70, 8).Value = 7399.4736
$ws.Cells.Item(70, 9).Value = 7949.1665
$ws.Cells.Item(70, 10).Value = 6457.143
$ws.Cells.Item(70, 11).Value = 7949.1665
$ws.Cells.Item(70, 12).Value = 6457.143
$ws.Cells.Item(70, 13).Value = -7679.1665
$ws.Cells.Item(70, 14).Value = -6997.143
$ws.Cells.Item(73, 8).Value = 7399.4736
$ws.Cells.Item(73, 9).Value = 7949.1665
$ws.Cells.Item(73, 10).Value = 6457.143
$ws.Cells.Item(73, 11).Value = 7949.1665
$ws.Cells.Item(73, 12).Value = 6457.143
$ws.Cells.Item(73, 13).Value = -7013.1665
$ws.Cells.Item(73, 14).Value = -8329.143
$ws.Cells.Item(122, 8).Value = 795552.94
$ws.Cells.Item(122, 9).Value = 1112703.2
$ws.Cells.Item(122, 10).Value = 2677
$ws.Cells.Item(122, 11).Value = 3338109.6
$ws.Cells.Item(122, 12).Value = 8031
$ws.Cells.Item(122, 13).Value = -3335659.6
$ws.Cells.Item(122, 14).Value = -12931
$ws.Cells.Item(132, 8).Value = 3109.1345
$ws.Cells.Item(132, 9).Value = 3013.7354
$ws.Cells.Item(132, 10).Value = 3289.3333
$ws.Cells.Item(132, 11).Value = 9041.2062
$ws.Cells.Item(132, 12).Value = 9867.999899999999
$ws.Cells.Item(132, 13).Value = -6511.206200000001
$ws.Cells.Item(132, 14).Value = -14927.9999
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 14).Value = $null
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 14).Value = $null
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3414.923
$ws.Cells.Item(40, 9).Value = 3197
$ws.Cells.Item(40, 10).Value = 3454.5454
$ws.Cells.Item(40, 11).Value = 3197
$ws.Cells.Item(40, 12).Value = 3454.5454
$ws.Cells.Item(40, 13).Value = -3061
$ws.Cells.Item(40, 14).Value = -3726.5454
$ws.Cells.Item(61, 8).Value = 3161.5334
$ws.Cells.Item(61, 9).Value = 5915
$ws.Cells.Item(61, 10).Value = 1325.8889
$ws.Cells.Item(61, 11).Value = 5915
$ws.Cells.Item(61, 12).Value = 1325.8889
$ws.Cells.Item(61, 13).Value = -5713
$ws.Cells.Item(61, 14).Value = -1729.8889
$ws.Cells.Item(113, 8).Value = 3161.5334
$ws.Cells.Item(113, 9).Value = 5915
$ws.Cells.Item(113, 10).Value = 1325.8889
$ws.Cells.Item(113, 11).Value = 5915
$ws.Cells.Item(113, 12).Value = 1325.8889
$ws.Cells.Item(113, 13).Value = -3745
$ws.Cells.Item(113, 14).Value = -5665.8889

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 50027500
$ws.Cells.Item(2, 9).Value = 57164284
$ws.Cells.Item(2, 11).Value = 57164284
$ws.Cells.Item(2, 13).Value = -57164172
$ws.Cells.Item(4, 8).Value = 21026.25
$ws.Cells.Item(4, 9).Value = 51
$ws.Cells.Item(4, 10).Value = 42001.5
$ws.Cells.Item(4, 11).Value = 51
$ws.Cells.Item(4, 12).Value = 42001.5
$ws.Cells.Item(4, 13).Value = 62
$ws.Cells.Item(4, 14).Value = -42227.5
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 14).Value = $null
$ws.Cells.Item(132, 8).Value = 1973.8229
$ws.Cells.Item(132, 9).Value = 2067.7778
$ws.Cells.Item(132, 10).Value = 1691.9584
$ws.Cells.Item(132, 11).Value = 6203.3334
$ws.Cells.Item(132, 12).Value = 5075.8752
$ws.Cells.Item(132, 13).Value = -3673.3334
$ws.Cells.Item(132, 14).Value = -10135.8752
